$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.003837439598571412
$ws.Range("E2").Value = 0.3718167378372522
$ws.Range("G2").Value = 0.2494892361375043
$ws.Range("I2").Value = 0.3694142004502854
$ws.Range("L2").Value = 0.5971552000000001
$ws.Range("M2").Value = 0.0822565
$ws.Range("N2").Value = 12.82009457445577
$ws.Range("O2").Value = 3.536987855065788

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.004836069939105167
$ws.Range("B2").Value = 0.04878002426392428
$ws.Range("E2").Value = 0.2223185247529475
$ws.Range("I2").Value = 0.4164466220504326
$ws.Range("L2").Value = 0.1199827757707056
$ws.Range("M2").Value = 0.04735360000000001
$ws.Range("N2").Value = 5.453610051398019
$ws.Range("O2").Value = 2.355488407130716

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.08761853486597547
$ws.Range("B2").Value = 0.02806921357332795
$ws.Range("E2").Value = 0.1716889961803926
$ws.Range("I2").Value = 0.4658448658569125
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04793623333333344
$ws.Range("N2").Value = 7.992614123745049
$ws.Range("O2").Value = 4.912910293249987
